$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = -0.213011757062836
$ws.Range("J4").Value = 0.4629671583968761
$ws.Range("K4").Value = 0.4214911133908249
$ws.Range("L4").Value = 2.799257955029619
